$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source Information")

$i2 = @(
  '1. The command "Table.TransformColumnTypes" transforms the data types of certain columns in a source table to text.',
  '2. "Promoted Headers" makes the first row of a table the header row and changes all data below it into a single type of data.',
  "3. This changes the data type of the 'Nominal GDP (INR,?)' and 'Zonal councils' columns from their original type to text.",
  '4. Renamed Columns changes the column names of "Changed Type1" to "Zone" and "Nominal GDP(INR, ?)" to "Nominal GDP".',
  "5. This replaces any question marks in the 'Nominal GDP' column with blank values.",
  '6. This sentence instructs the computer to replace the phrase "lakh crore" with nothing in the "Nominal GDP" column of a table.',
  '7. This sentence changes the data type of the column "Nominal GDP" to a type of Currency.'
) -join "`n`n"
$i2 = $i2 + "`n`n"

$i3 = @(
  "1. This changes the data types of the columns in the table 'GDP20012010_Table' to text and any.",
  '2. "Unpivoting Other Columns" means converting multiple columns into two, one for attributes and one for their corresponding values.',
  '3. "Appended Query" means combining two tables, "Unpivoted Other Columns" and "GDP (2011-2014)", into one table.',
  '4. The "Attribute" column has been renamed to "Year" in the Appended Query table.',
  '5. Table.TransformColumnTypes changes the data type of the "Year" column from its previous data type to "date".',
  '6. Table "Changed Type1" has been renamed to "Renamed Columns1" and the column "Value" has been changed to "GDP".',
  '7. "Changed Type2" is transforming the "GDP" column into a currency type.',
  '8. "Replaced Errors" replaces any errors in the "Changed Type2" table with 0 for the GDP column.',
  '9. Merged Queries will combine the information from two tables - "Replaced Errors" and "GDP Growth Rate (2002-2014)" - by joining the matching entries of "State/union territory" and "Year" and display all entries from the "Replaced Errors" table plus the related matching entries',
  '10. This operation takes the GDP Growth Rate (2002-2014) information from one table and expands it into three separate columns containing the State/Union territory, Year, and Growth Rate(%) information.',
  '11. This statement removes the columns labeled "State/union territory" and "Year" from the table "Expanded GDP Growth Rate (2002-2014)".',
  '12. RenameColumns2 changes the names of two columns from "GDP Growth Rate (2002-2014).Growth Rate(%)" to "Growth Rate(%)" and from "State/union territory" to "State".',
  "13. This sentence replaces any null values in the 'Growth Rate(%)' column with 0s."
) -join "`n`n"
$i3 = $i3 + "`n`n"

$i4 = @(
  '1. "Changed Type" means that the table GrowthRate20022010_Table was modified so that the columns "2002" to "2010" are treated as any type, and the column "State/union territory" is treated as text.',
  "2. This command removes any rows with errors from a table that had its datatype changed, and specifically the row containing the value 2002.",
  '3. Unpivoted Other Columns is a process of turning data from multiple columns into two columns, "Attribute" and "Value".',
  "4. Appending a Query combines the results from two separate tables into one.",
  '5. "Renamed Columns" is the new table name that replaces "Appended Query" after its column named "Attribute" is renamed to "Year".',
  '6. This means the columns in the Renamed Columns table have been changed so that the column "Year" is now a date type.',
  "7. The table's Value column was renamed to Growth Rate(%) using Table.RenameColumns.",
  '8. This statement changes the data type of the "Growth Rate (%)" column from the previous step to a percentage type.'
) -join "`n`n"
$i4 = $i4 + "`n`n"

$i5 = @(
  "1. The Changed Type command changes the data types of each column in the GrowthRate20112014_Table table to text, any, or any.",
  '2. Removed any rows with errors in the column "2011" from the table which was previously changed type.',
  '3. This code removes any rows with errors from the "Removed Errors" table, specifically from the column labeled "2013".',
  "4. Unpivoting Other Columns rearranges the columns in the Removed Errors1 table, putting all the State/union territory values in one column and the related attributes in another."
) -join "`n`n"
$i5 = $i5 + "`n`n"

$i6 = @(
  "1. Promoted Headers is a function that takes table columns and combines them into one header row.",
  '2. "Changed Type1" changed the data type of the columns "State/union territory", "2011", "2012", "2013", "2014", and "Column6" to text, any, and any format, respectively.',
  '3. "Removed Columns" means to remove the column named "Column6" from the table "Changed Type1".',
  "4. Unpivoting other columns reorganizes the data so that each column is converted into multiple rows with additional columns indicating the original column values."
) -join "`n`n"
$i6 = $i6 + "`n`n"

$ws.Range("I2").Value = $i2
$ws.Range("I3").Value = $i3
$ws.Range("I4").Value = $i4
$ws.Range("I5").Value = $i5
$ws.Range("I6").Value = $i6
